$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Actual" temperatures that had just become available for the
# two most recent forecast rows (May/June). The adjacent "Diff (fc-act)"
# column recalculates automatically from its existing D-E formula.
$ws.Range("E13").Value = 20.3
$ws.Range("E14").Value = 19.600000000000001

# Leave the selection where the user's cursor ended up after the entry.
$ws.Range("E15").Select() | Out-Null
